# DOMA-4452: add "Место установки счетчика" (meter installation place) column
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column S inherits the same look (style/border/fill/font) as column R, which
# is the last column of the existing header/data table.
$ws.Range("R1:R11").Copy()
$ws.Range("S1:S11").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("S:S").ColumnWidth = $ws.Range("R:R").ColumnWidth

# Header
$ws.Range("S1").Value = "Место установки счетчика"

# Sample data - alternating "Kitchen" / "Bathroom" meter locations
$ws.Range("S2").Value = "Кухня"
$ws.Range("S3").Value = "Сан. узел"
$ws.Range("S4").Value = "Кухня"
$ws.Range("S5").Value = "Сан. узел"
$ws.Range("S6").Value = "Кухня"
$ws.Range("S7").Value = "Сан. узел"
$ws.Range("S8").Value = "Кухня"
$ws.Range("S9").Value = "Сан. узел"
$ws.Range("S10").Value = "Кухня"
$ws.Range("S11").Value = "Сан. узел"
